# Add team record (Wins/Losses/Ties) columns to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (AC1)
# onto the three new header cells so they match the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels.
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Team record values for every player row (2-52).
$ws.Range("AD2:AD52").Value2 = 89
$ws.Range("AE2:AE52").Value2 = 73
$ws.Range("AF2:AF52").Value2 = 0
